# Update scripts with new tpm
# - Recompute Ligand/Receptor/Edge specificity columns (I:T) for the two
#   remaining FAPs-sending rows now that the MuSCs-sending rows are gone.
# - Remove the two rows where "MuSCs" was the sending cluster (old rows 4 & 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (FAPs -> Rspo2/Lgr6 -> FAPs): refreshed specificity values ---
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.05147733333333334
$ws.Range("O2").Value = 0.887188413789934
$ws.Range("P2").Value = 0.8871884137899338
$ws.Range("Q2").Value = 0.07859821787733334
$ws.Range("R2").Value = 0.7073839608960001
$ws.Range("S2").Value = 0.887188413789934
$ws.Range("T2").Value = 0.8871884137899338

# --- Row 3 (FAPs -> Rspo2/Lgr6 -> MuSCs): refreshed specificity values ---
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.006545666666666668
$ws.Range("N3").Value = 0.019637
$ws.Range("O3").Value = 0.1128115862100661
$ws.Range("P3").Value = 0.1128115862100661
$ws.Range("Q3").Value = 0.009994257695666668
$ws.Range("R3").Value = 0.08994831926100001
$ws.Range("S3").Value = 0.1128115862100661
$ws.Range("T3").Value = 0.1128115862100661

# --- Remove old rows 4 and 5 (MuSCs as sending cluster) ---
$ws.Rows(4).EntireRow.Delete()
$ws.Rows(4).EntireRow.Delete()
